$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the header formatting
# (bold font, border, centered alignment) from the existing header cell E1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Populate the time_taken column for every data row (2-72) with the
# recorded timestamp values, written as plain text.
$timeTakenValues = @(
    "2021-10-05 10:51:54.062253",
    "2021-10-05 10:51:54.062264",
    "2021-10-05 10:51:54.062267",
    "2021-10-05 10:51:54.062270",
    "2021-10-05 10:51:54.062273",
    "2021-10-05 10:51:54.062275",
    "2021-10-05 10:51:54.062278",
    "2021-10-05 10:51:54.062281",
    "2021-10-05 10:51:54.062283",
    "2021-10-05 10:51:54.062286",
    "2021-10-05 10:51:54.062288",
    "2021-10-05 10:51:54.062291",
    "2021-10-05 10:51:54.062294",
    "2021-10-05 10:51:54.062296",
    "2021-10-05 10:51:54.062299",
    "2021-10-05 10:51:54.062301",
    "2021-10-05 10:51:54.062304",
    "2021-10-05 10:51:54.062307",
    "2021-10-05 10:51:54.062309",
    "2021-10-05 10:51:54.062312",
    "2021-10-05 10:51:54.062314",
    "2021-10-05 10:51:54.062317",
    "2021-10-05 10:51:54.062320",
    "2021-10-05 10:51:54.062322",
    "2021-10-05 10:51:54.062325",
    "2021-10-05 10:51:54.062327",
    "2021-10-05 10:51:54.062330",
    "2021-10-05 10:51:54.062332",
    "2021-10-05 10:51:54.062335",
    "2021-10-05 10:51:54.062340",
    "2021-10-05 10:51:54.062342",
    "2021-10-05 10:51:54.062345",
    "2021-10-05 10:51:54.062348",
    "2021-10-05 10:51:54.062350",
    "2021-10-05 10:51:54.062353",
    "2021-10-05 10:51:54.062356",
    "2021-10-05 10:51:54.062358",
    "2021-10-05 10:51:54.062361",
    "2021-10-05 10:51:54.062363",
    "2021-10-05 10:51:54.062366",
    "2021-10-05 10:51:54.062369",
    "2021-10-05 10:51:54.062372",
    "2021-10-05 10:51:54.062374",
    "2021-10-05 10:51:54.062377",
    "2021-10-05 10:51:54.062379",
    "2021-10-05 10:51:54.062382",
    "2021-10-05 10:51:54.062384",
    "2021-10-05 10:51:54.062387",
    "2021-10-05 10:51:54.062390",
    "2021-10-05 10:51:54.062392",
    "2021-10-05 10:51:54.062395",
    "2021-10-05 10:51:54.062397",
    "2021-10-05 10:51:54.062400",
    "2021-10-05 10:51:54.062403",
    "2021-10-05 10:51:54.062406",
    "2021-10-05 10:51:54.062408",
    "2021-10-05 10:51:54.062410",
    "2021-10-05 10:51:54.062413",
    "2021-10-05 10:51:54.062415",
    "2021-10-05 10:51:54.062418",
    "2021-10-05 10:51:54.062420",
    "2021-10-05 10:51:54.062423",
    "2021-10-05 10:51:54.062426",
    "2021-10-05 10:51:54.062428",
    "2021-10-05 10:51:54.062432",
    "2021-10-05 10:51:54.062435",
    "2021-10-05 10:51:54.062437",
    "2021-10-05 10:51:54.062440",
    "2021-10-05 10:51:54.062442",
    "2021-10-05 10:51:54.062445",
    "2021-10-05 10:51:54.062447"
)

for ($i = 0; $i -lt $timeTakenValues.Length; $i++) {
    $targetRow = $i + 2
    $ws.Cells.Item($targetRow, 6).Value = $timeTakenValues[$i]
}
